# ============================================================================
# Applies the "madrugada" revision pass to "LEIAM COM ATENÇÃO.docx":
#   - wraps a handful of foreign/slang words (etc, to, hash, DDoS, IDEs,
#     JavaScript, main) in <w:proofErr spellStart/spellEnd> pairs, which
#     forces Word to split the host run around them (mirrors what happens
#     when Word's background spell-checker flags a word as it is typed)
#   - appends new trailing sentences to the "Decidi também..." paragraph
#   - appends a brand-new paragraph with the "Realizei pequenos testes..."
#     status update
#
# NOTE: this interpreter's argument parser chokes on `<`/`>` characters when
# they show up in a *named* parameter value (e.g. `-Foo $bar` where $bar
# contains XML) - so every helper below is called positionally.
# ============================================================================

$d = $word.ActiveDocument

$wordmlNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Set-ParagraphInnerXml($Paragraph, $InnerXml) {
    $pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
           '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData>' +
           '<w:document ' + $wordmlNs + '><w:body>' + $InnerXml + '</w:body></w:document>' +
           '</pkg:xmlData></pkg:part></pkg:package>'
    $Paragraph.Range.InsertXML($pkg)
}

function New-Run($Text, $Preserve) {
    if ($Preserve) {
        return '<w:r><w:t xml:space="preserve">' + $Text + '</w:t></w:r>'
    } else {
        return '<w:r><w:t>' + $Text + '</w:t></w:r>'
    }
}

function New-ProofedRun($Text) {
    return '<w:proofErr w:type="spellStart"/>' + (New-Run $Text $false) + '<w:proofErr w:type="spellEnd"/>'
}

# ---------------------------------------------------------------------------
# Paragraph 3: "Galera seguinte, ..." -> split around "etc" and "to"
# ---------------------------------------------------------------------------
$p3 = $d.Paragraphs(3)
$p3inner = '<w:p>' +
    (New-Run 'Galera seguinte, rota final do projeto e eu sei que somos pessoas ocupadas e nem todos tiveram o mesmo tempo produzindo ' $true) +
    (New-ProofedRun 'etc') +
    (New-Run ', acima de tudo somos uma equipe e sei que cada um contribuiu como pôde. Porém preciso que vocês acessem esse arquivo aqui e deixem anotado as ideias que tiveram e suas dificuldades na hora de desenvolver. Por exemplo, eu, ' $true) +
    (New-ProofedRun 'to') +
    (New-Run ' com uma certa dificuldade em identificar os erros e achar as soluções para isso. Anotem o que puderem que eu passo corrigindo e formatando para ABNT e reviso depois com o Matheus. Tamo junto galera!!' $true) +
    '</w:p>'
Set-ParagraphInnerXml $p3 $p3inner

# ---------------------------------------------------------------------------
# Paragraph 5: "Matheus: ..." -> split around both occurrences of "hash"
# ---------------------------------------------------------------------------
$p5 = $d.Paragraphs(5)
$p5inner = '<w:p>' +
    (New-Run 'Matheus: Obtive problemas em criar uma conta administrativa e normal, pois por segurança as senhas são salvas em ' $true) +
    (New-ProofedRun 'hash') +
    (New-Run ', logo tive que pesquisar como que transformava uma senha de caracteres em ' $true) +
    (New-ProofedRun 'hash') +
    '</w:p>'
Set-ParagraphInnerXml $p5 $p5inner

# ---------------------------------------------------------------------------
# Paragraph 7: "Marcos: ..." -> split around "DDoS"
# ---------------------------------------------------------------------------
$p7 = $d.Paragraphs(7)
$p7inner = '<w:p>' +
    (New-Run 'Marcos: Encontrei difi' $false) +
    (New-Run 'culdade em ajustar e otimizar o código. Tinha algumas pequenas coisas há melhorar e algumas pequenas correções a fazer no script SQL para que os testes fossem realizados de forma consistente.' $false) +
    (New-Run ' A consistência e construção do código o torna quase perfeito e de difícil acesso ' $true) +
    (New-Run 'ilegal, surge a ideia de fazer um CAPTCHA para garantir que o site não sofra um ataque ' $true) +
    (New-ProofedRun 'DDoS') +
    (New-Run ', porém não se mostra necessário devido ser um site de cadastro' $false) +
    (New-Run ' de doadores de sangue e agendamento de doação de sangue.' $true) +
    '</w:p>'
Set-ParagraphInnerXml $p7 $p7inner

# ---------------------------------------------------------------------------
# Paragraph 8: "Decidi também ..." -> split around "IDEs" and "JavaScript",
# plus three new trailing sentences after the existing <w:br/> line.
# ---------------------------------------------------------------------------
$p8 = $d.Paragraphs(8)
$p8inner = '<w:p>' +
    (New-Run 'Decidi também separar a parte de estilização em dois arquivos diferentes' $false) +
    (New-Run ',' $false) +
    (New-Run ' para que o site ficasse mais otimizado, também facilitando a visualização do código em ' $true) +
    (New-ProofedRun 'IDEs') +
    (New-Run '. Segreguei parte do código deixando somente a parte funcional em evidência nos arquivos HTML e PHP' $false) +
    (New-Run '. A razão da segregação dessa parte do código tem relação com a otimização do site, a parte de ' $true) +
    (New-ProofedRun 'JavaScript') +
    (New-Run ' não precisou ser separada em arquivos diferentes' $true) +
    (New-Run ', parte bem leve do código.' $false) +
    '<w:r><w:br/><w:t>Encontrei certa dificuldade em fazer essa segregação para garantir de que não houvesse erros</w:t></w:r>' +
    (New-Run ' e posso afirmar que foi um sucesso.' $true) +
    (New-Run ' Levei horas fazendo e fui até a falha, madrugada porém valeu a pena cada segundo e hora fazendo essa revisão' $true) +
    (New-Run '.' $false) +
    '</w:p>'
Set-ParagraphInnerXml $p8 $p8inner

# ---------------------------------------------------------------------------
# Brand-new paragraph right after paragraph 8 ("Realizei pequenos testes...")
# ---------------------------------------------------------------------------
$p8 = $d.Paragraphs(8)
$p8.Range.InsertParagraphAfter()
$pNew = $d.Paragraphs(9)
$pNewInner = '<w:p>' +
    (New-Run 'Realizei pequenos testes de campo e tudo ocorreu bem, ' $true) +
    (New-Run 'cada um com seu cada qual funcionando perfeitamente. Fiz um merge das melhorias e inclui algumas coisas no banco de dados' $false) +
    (New-Run ' para a realização de alguns testes caso necessário em apresentação na quinta feira dia 27' $true) +
    (New-Run '. ' $true) +
    (New-Run 'Creio que só slides não serão ' $true) +
    (New-Run 'suficientes' $false) +
    (New-Run ' então justamente pra garantir, deixei na ' $true) +
    (New-ProofedRun 'main') +
    (New-Run ' para que possamos fazer o famoso quadro: "Quem Sabe Faz Ao Vivo"!' $true) +
    '</w:p>'
Set-ParagraphInnerXml $pNew $pNewInner

Write-Host "Done."
